$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.Formula = "'29.225.80"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +0.02%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.Formula = "'1.848.64"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -0.63%  "
$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.Formula = "'0.9987"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.Formula = "'246.40"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +2.02%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.Formula = "'0.7001"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -0.91%  "
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.Formula = "'0.9993"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.05%  "
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.Formula = "'0.07731"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -1.52%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.Formula = "'23.59"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -1.10%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.Formula = "'0.07826"
$cell.Style = $origStyle
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.Formula = "'93.30"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +0.72%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.Formula = "'1.845.52"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -0.55%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.Formula = "'5.131"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +0.14%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.Formula = "'0.6875"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -0.27%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.Formula = "'6.638"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  -1.52%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.Formula = "'29.196.95"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +0.00%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.Formula = "'241.44"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -3.57%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.Formula = "'2.083.47"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -0.62%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.Formula = "'12.78"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -0.95%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.Formula = "'0.9992"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  -0.05%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.Formula = "'7.525"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -0.02%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.Formula = "'0.1519"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -1.16%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.Formula = "'158.99"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -0.77%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.Formula = "'8.837"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  -1.31%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.Formula = "'4.239"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -0.94%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.Formula = "'4.194"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -1.34%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.Formula = "'1.197"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -0.58%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.Formula = "'0.05127"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -1.58%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.Formula = "'0.7926"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +4.31%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.Formula = "'1.873"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  -2.13%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.Formula = "'2.690"
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.Formula = "'1.312.65"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +6.84%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.Formula = "'0.01873"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +0.62%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.Formula = "'2.710"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -0.42%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.Formula = "'0.9490"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +5.64%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.Formula = "'6.080"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +6.60%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.Formula = "'107.74"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E44").Value = "  +0.00%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.Formula = "'9.739"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("E46").Value = "  -0.91%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.Formula = "'1.985.06"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -0.61%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.Formula = "'0.5177"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -1.64%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.Formula = "'1.769"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +0.54%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.Formula = "'7.012"
$cell.Style = $origStyle
